$d = $word.ActiveDocument

# --- Locate the trailing paragraph that currently holds only a tab and the
# --- hidden "_GoBack" bookmark (the last paragraph in the body). ---
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$paraStart = $lastPara.Range.Start

# Remove the existing _GoBack bookmark; we will re-create it in the exact
# spot dictated by the target markup (immediately after the new run).
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete() | Out-Null

# Replace the tab character (the sole content of that paragraph) with a run
# that keeps the tab and appends the new sentence, followed by the restored
# bookmark - matching the target paragraph structure exactly.
$tabRange = $d.Range($paraStart, $paraStart + 1)
$xmlFirst = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t>In addition to unread, nots and zeros can also be done in reverse</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$tabRange.InsertXML($xmlFirst) | Out-Null

# --- Append three brand-new paragraphs after that paragraph, right before
# --- the section break, reproducing the remaining part of the diff. ---
$insertAt = $d.Paragraphs.Item($lastParaIndex).Range.End
$tail = $d.Range($insertAt, $insertAt)
$xmlTail = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:tab/></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Types of errors:</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:tab/><w:t>Using the TEMP in a map when it doesn’t exist</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$tail.InsertXML($xmlTail) | Out-Null
